$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 55556484
$ws.Range("I18").Value = 55556484
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 55556484
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -55556200
$ws.Range("H32").Value = 2899.75
$ws.Range("I32").Value = 3566.1667
$ws.Range("J32").Value = 2677.611
$ws.Range("K32").Value = 3566.1667
$ws.Range("L32").Value = 2677.611
$ws.Range("M32").Value = -3240.1667
$ws.Range("N32").Value = -3329.611
$ws.Range("H43").Value = 1930139.2
$ws.Range("I43").Value = 2572186.2
$ws.Range("J43").Value = 3998
$ws.Range("K43").Value = 2572186.2
$ws.Range("L43").Value = 3998
$ws.Range("M43").Value = -2572117.2
$ws.Range("N43").Value = -4136
$ws.Range("H98").Value = 1702.92
$ws.Range("I98").Value = 1726.15
$ws.Range("J98").Value = 1610
$ws.Range("K98").Value = 1726.15
$ws.Range("L98").Value = 1610
$ws.Range("M98").Value = -228.1500000000001
$ws.Range("N98").Value = -4606
$ws.Range("H116").Value = 2088335
$ws.Range("I116").Value = 3792581
$ws.Range("J116").Value = 5367.5557
$ws.Range("K116").Value = 3792581
$ws.Range("L116").Value = 5367.5557
$ws.Range("M116").Value = -3789139
$ws.Range("N116").Value = -12251.5557
$ws.Range("H122").Value = 1702.92
$ws.Range("I122").Value = 1726.15
$ws.Range("J122").Value = 1610
$ws.Range("K122").Value = 5178.450000000001
$ws.Range("L122").Value = 4830
$ws.Range("M122").Value = -2728.450000000001
$ws.Range("N122").Value = -9730
$ws.Range("H132").Value = 9920.781000000001
$ws.Range("I132").Value = 4776.0513
$ws.Range("J132").Value = 13440.859
$ws.Range("K132").Value = 14328.1539
$ws.Range("L132").Value = 40322.577
$ws.Range("M132").Value = -11798.1539
$ws.Range("N132").Value = -45382.577
$ws.Range("H140").Value = 71882.5
$ws.Range("I140").Value = 73333.336
$ws.Range("J140").Value = 69706.25
$ws.Range("K140").Value = 73333.336
$ws.Range("L140").Value = 69706.25
$ws.Range("M140").Value = -68153.336
$ws.Range("N140").Value = -80066.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 700251.8
$ws.Range("I2").Value = 1061810
$ws.Range("J2").Value = 57481.668
$ws.Range("K2").Value = 1061810
$ws.Range("L2").Value = 57481.668
$ws.Range("M2").Value = -1061697
$ws.Range("N2").Value = -57707.668
$ws.Range("H45").Value = 1740.5
$ws.Range("I45").Value = 1602.2858
$ws.Range("J45").Value = 2063
$ws.Range("K45").Value = 1602.2858
$ws.Range("L45").Value = 2063
$ws.Range("M45").Value = -1225.2858
$ws.Range("N45").Value = -2817
$ws.Range("H97").Value = 715.1875
$ws.Range("I97").Value = 723.7143
$ws.Range("J97").Value = 655.5
$ws.Range("K97").Value = 723.7143
$ws.Range("L97").Value = 655.5
$ws.Range("M97").Value = -227.7143
$ws.Range("N97").Value = -1647.5
$ws.Range("H116").Value = 700251.8
$ws.Range("I116").Value = 1061810
$ws.Range("J116").Value = 57481.668
$ws.Range("K116").Value = 1061810
$ws.Range("L116").Value = 57481.668
$ws.Range("M116").Value = -1059516
$ws.Range("N116").Value = -62069.668
$ws.Range("H122").Value = 3597.8333
$ws.Range("I122").Value = 1966.0834
$ws.Range("J122").Value = 10124.833
$ws.Range("K122").Value = 5898.2502
$ws.Range("L122").Value = 30374.499
$ws.Range("M122").Value = -3448.2502
$ws.Range("N122").Value = -35274.499
$ws.Range("H140").Value = 110148.836
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 110148.836
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 110148.836
$ws.Range("N140").Value = -120508.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 700251.8
$ws.Range("I3").Value = 1061810
$ws.Range("J3").Value = 57481.668
$ws.Range("K3").Value = 1061810
$ws.Range("L3").Value = 57481.668
$ws.Range("M3").Value = -1061696
$ws.Range("N3").Value = -57709.668
$ws.Range("H94").Value = 856981.0600000001
$ws.Range("I94").Value = 1712837.2
$ws.Range("J94").Value = 1124.875
$ws.Range("K94").Value = 1712837.2
$ws.Range("L94").Value = 1124.875
$ws.Range("M94").Value = -1712386.2
$ws.Range("N94").Value = -2026.875
$ws.Range("H99").Value = 1159425.5
$ws.Range("I99").Value = 1227450.5
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 1227450.5
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1225952.5
$ws.Range("N99").Value = -5996
$ws.Range("H134").Value = 2956.6875
$ws.Range("I134").Value = 2820.4666
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 8461.399800000001
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -5926.399800000001
$ws.Range("N134").Value = -20070
$ws.Range("H140").Value = 138403.36
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 138403.36
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 138403.36
$ws.Range("N140").Value = -148763.36

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3877.1155
$ws.Range("I31").Value = 3197.0278
$ws.Range("J31").Value = 5407.3125
$ws.Range("K31").Value = 3197.0278
$ws.Range("L31").Value = 5407.3125
$ws.Range("M31").Value = -2902.0278
$ws.Range("N31").Value = -5997.3125
$ws.Range("H34").Value = 3877.1155
$ws.Range("I34").Value = 3197.0278
$ws.Range("J34").Value = 5407.3125
$ws.Range("K34").Value = 3197.0278
$ws.Range("L34").Value = 5407.3125
$ws.Range("M34").Value = -2995.0278
$ws.Range("N34").Value = -5811.3125
$ws.Range("H62").Value = 47666.555
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 59999.855
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 59999.855
$ws.Range("M62").Value = -3876
$ws.Range("N62").Value = -61247.855
$ws.Range("H65").Value = 47666.555
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 59999.855
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 299999.275
$ws.Range("M65").Value = -19380
$ws.Range("N65").Value = -306239.275
$ws.Range("H99").Value = 23758.375
$ws.Range("I99").Value = 40022.332
$ws.Range("J99").Value = 14000
$ws.Range("K99").Value = 40022.332
$ws.Range("L99").Value = 14000
$ws.Range("M99").Value = -38524.332
$ws.Range("N99").Value = -16996
$ws.Range("H105").Value = 1624249.9
$ws.Range("I105").Value = 2066792.9
$ws.Range("J105").Value = 1592.3334
$ws.Range("K105").Value = 2066792.9
$ws.Range("L105").Value = 1592.3334
$ws.Range("M105").Value = -2065045.9
$ws.Range("N105").Value = -5086.3334
$ws.Range("H109").Value = 59091.668
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 59091.668
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 59091.668
$ws.Range("N109").Value = -61171.668
$ws.Range("H126").Value = 23758.375
$ws.Range("I126").Value = 40022.332
$ws.Range("J126").Value = 14000
$ws.Range("K126").Value = 120066.996
$ws.Range("L126").Value = 42000
$ws.Range("M126").Value = -117596.996
$ws.Range("N126").Value = -46940
$ws.Range("H141").Value = 83335.55499999999
$ws.Range("I141").Value = 41999.5
$ws.Range("J141").Value = 88502.56
$ws.Range("K141").Value = 41999.5
$ws.Range("L141").Value = 88502.56
$ws.Range("M141").Value = -36819.5
$ws.Range("N141").Value = -98862.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 136.76471
$ws.Range("I2").Value = 80.55556
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 483.33336
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -370.33336
$ws.Range("N2").Value = -1426
$ws.Range("H7").Value = 299
$ws.Range("I7").Value = 299
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 897
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -785
$ws.Range("H23").Value = 71428860
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 71428860
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 214286580
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -214287050
$ws.Range("H34").Value = 2943.3333
$ws.Range("I34").Value = 3100
$ws.Range("J34").Value = 2818
$ws.Range("K34").Value = 9300
$ws.Range("L34").Value = 8454
$ws.Range("M34").Value = -9216
$ws.Range("N34").Value = -8622
$ws.Range("H55").Value = 4950
$ws.Range("I55").Value = 900
$ws.Range("J55").Value = 9000
$ws.Range("K55").Value = 2700
$ws.Range("L55").Value = 27000
$ws.Range("M55").Value = -2523
$ws.Range("N55").Value = -27354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 12909.777
$ws.Range("I113").Value = 11747.25
$ws.Range("J113").Value = 13839.8
$ws.Range("K113").Value = 11747.25
$ws.Range("L113").Value = 13839.8
$ws.Range("M113").Value = -9577.25
$ws.Range("N113").Value = -18179.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6456.5
$ws.Range("I132").Value = 4869.1665
$ws.Range("J132").Value = 7647
$ws.Range("K132").Value = 14607.4995
$ws.Range("L132").Value = 22941
$ws.Range("M132").Value = -12077.4995
$ws.Range("N132").Value = -28001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 11000
$ws.Range("I22").Value = 11000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 11000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -10707
$ws.Range("H126").Value = 1403.2106
$ws.Range("I126").Value = 1308
$ws.Range("J126").Value = 1488.9
$ws.Range("K126").Value = 3924
$ws.Range("L126").Value = 4466.700000000001
$ws.Range("M126").Value = -1454
$ws.Range("N126").Value = -9406.700000000001
$ws.Range("H132").Value = 13163302
$ws.Range("I132").Value = 1674.6333
$ws.Range("J132").Value = 62519404
$ws.Range("K132").Value = 5023.8999
$ws.Range("L132").Value = 187558212
$ws.Range("M132").Value = -2493.8999
$ws.Range("N132").Value = -187563272
$ws.Range("H138").Value = 76243.2
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 76243.2
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 76243.2
$ws.Range("N138").Value = -86523.2
